$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Id" values in column A for rows 16-19 (Major Changes / customer Ids renumbered)
$ws.Range("A16").Value = 124
$ws.Range("A17").Value = 125
$ws.Range("A18").Value = 126
$ws.Range("A19").Value = 127

# Move the active selection to A22 (mirrors the user's last selection when saving)
$ws.Range("A22").Select()
